$wb = $excel.ActiveWorkbook

# The two sheets "展览" and "全部类型" contain identical data tables and both
# need the same updates to column F ("想去人数").
$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new F-column value
$updates = @{
    2  = 821
    4  = 1145
    6  = 12334
    11 = 1129
    12 = 908
    13 = 13608
    14 = 13790
    15 = 41
    22 = 51
    23 = 4885
    24 = 214
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
